# Generate Report for Archive
# 1) Update shared "Ready for handoff" status text to "In Translation"
#    across all sheets that reference it.
# 2) Narrow the "zh-cn"/"de-de" (Overview) and "Status" (per-locale) columns
#    from 17.2159881591797 to 13.4101845877511 characters wide.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# -- Text update: replace "Ready for handoff" with "In Translation" --
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# -- Column width update --
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511
$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
